# TC_63799 - add "DC Unit Loading Details" column (E1:E3) to the two
# loop sheets and flip which sheet/selection is active, matching the
# "Updated 19 test cases with loop loading details method" commit.

$wb = $excel.ActiveWorkbook

$wsAdd   = $wb.Worksheets.Item("Add Devices Loop A")
$wsOther = $wb.Worksheets.Item("Other Devices Loop A")

# ---- "Add Devices Loop A" sheet: insert new header/body cells E1:E3 ----
# E1 picks up the bold/blue header style already used by row 7 (A7:N7).
$wsAdd.Range("A7").Copy()
$wsAdd.Range("E1").PasteSpecial(-4122)
$wsAdd.Range("E1").Value = "DC Unit Loading Details Name"

# E2 / E3 pick up the shaded, left-aligned, wrapped body style used by D8.
$wsAdd.Range("D8").Copy()
$wsAdd.Range("E2").PasteSpecial(-4122)
$wsAdd.Range("E2").Value = "Current (DC Units)"

$wsAdd.Range("D8").Copy()
$wsAdd.Range("E3").PasteSpecial(-4122)
$wsAdd.Range("E3").Value = "Current (worst case)"

# ---- "Other Devices Loop A" sheet: same new column ----
$wsOther.Range("A7").Copy()
$wsOther.Range("E1").PasteSpecial(-4122)
$wsOther.Range("E1").Value = "DC Unit Loading Details Name"

$wsOther.Range("D8").Copy()
$wsOther.Range("E2").PasteSpecial(-4122)
$wsOther.Range("E2").Value = "Current (DC Units)"

$wsOther.Range("D8").Copy()
$wsOther.Range("E3").PasteSpecial(-4122)
$wsOther.Range("E3").Value = "Current (worst case)"

# ---- Selections / active tab flip ----
# Before: "Other Devices Loop A" tab selected, selection C9.
# After:  "Add Devices Loop A" tab selected, selection E1:E3 (anchor E1);
#         "Other Devices Loop A" keeps selection E2:E3 (anchor E2).
$wsOther.Range("E2:E3").Select() | Out-Null

$wsAdd.Activate() | Out-Null
$wsAdd.Range("E1:E3").Select() | Out-Null
